$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Delete entire row 151 ("Apply AdminLTE Theme..." task), shifting everything below up by one.
$ws.Rows.Item(151).Delete()

# Update the frozen-pane top-left cell and active selection to match final view state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 146
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E152").Select()
